$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text updates (row 1) ---
# A1: "订单编号" -> "*订单编号"
$ws.Range("A1").Value = "*订单编号"
# C1: "商家ID" -> "*商家ID"
$ws.Range("C1").Value = "*商家ID"
# D1: text stays "商家名称", but becomes centered + red font
$ws.Range("D1").Value = "商家名称"
# E1: "订单履约状态" -> "*订单履约状态"
$ws.Range("E1").Value = "*订单履约状态"

# --- Style updates ---
# E1 should end up with the same centered style already used by A1/C1 (font with family=3).
# Copy that format over instead of re-deriving a brand new font entry.
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# D1 should end up centered too, with a red font (new style). Start from the same
# centered/family=3 base font as C1/E1, then recolor it red so the engine derives the
# new font entry from that base rather than from the plain default font.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Font.Color = 255

$excel.CutCopyMode = 0

# --- Column width: column A width 13.44140625 -> 15 ---
$ws.Columns.Item(1).ColumnWidth = 14.285714285714286

# --- Selection change: B5 -> G19 ---
$ws.Range("G19").Select()
